$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: E1/F1/G1 renamed (A1-D1 and H1 stay the same) ---
$ws.Range("E1").Value = "numerator"
$ws.Range("F1").Value = "numerator_desc"
$ws.Range("G1").Value = "follow_up"

# --- Data rows 2-18: ids shift down a row (one new patient inserted, plus a
#     reordering of an existing row), and the follow_up/numerator flags are recomputed ---
$ws.Range("A2").Value = "11332-2024"
$ws.Range("B2").Value = "'11332"
$ws.Range("C2").Value = "'26332"
$ws.Range("E2").Value = $False
$ws.Range("F2").Value = "No screening recorded"
$ws.Range("G2").Value = $True
$ws.Range("H2").Value = $False

$ws.Range("A3").Value = "20202-2024"
$ws.Range("B3").Value = "'20202"
$ws.Range("C3").Value = "'30238"
$ws.Range("E3").Value = $False
$ws.Range("F3").Value = "No screening recorded"
$ws.Range("G3").Value = $False
$ws.Range("H3").Value = $False

$ws.Range("A4").Value = "26287-2024"
$ws.Range("B4").Value = "'26287"
$ws.Range("C4").Value = "'14630"
$ws.Range("E4").Value = $False
$ws.Range("F4").Value = "No screening recorded"
$ws.Range("G4").Value = $False
$ws.Range("H4").Value = $False

$ws.Range("A5").Value = "32851-2024"
$ws.Range("B5").Value = "'32851"
$ws.Range("C5").Value = "'58726"
$ws.Range("E5").Value = $False
$ws.Range("F5").Value = "No screening recorded"
$ws.Range("G5").Value = $True
$ws.Range("H5").Value = $False

$ws.Range("A6").Value = "34899-2024"
$ws.Range("B6").Value = "'34899"
$ws.Range("C6").Value = "'39533"
$ws.Range("E6").Value = $False
$ws.Range("F6").Value = "No screening recorded"
$ws.Range("G6").Value = $True
$ws.Range("H6").Value = $False

$ws.Range("A7").Value = "35383-2024"
$ws.Range("B7").Value = "'35383"
$ws.Range("C7").Value = "'62893"
$ws.Range("E7").Value = $False
$ws.Range("F7").Value = "No screening recorded"
$ws.Range("G7").Value = $False
$ws.Range("H7").Value = $False

$ws.Range("A8").Value = "36963-2024"
$ws.Range("B8").Value = "'36963"
$ws.Range("C8").Value = "'67362"
$ws.Range("E8").Value = $False
$ws.Range("F8").Value = "No screening recorded"
$ws.Range("G8").Value = $False
$ws.Range("H8").Value = $True

$ws.Range("A9").Value = "44221-2024"
$ws.Range("B9").Value = "'44221"
$ws.Range("C9").Value = "'64340"
$ws.Range("E9").Value = $False
$ws.Range("F9").Value = "No screening recorded"
$ws.Range("G9").Value = $False
$ws.Range("H9").Value = $False

$ws.Range("A10").Value = "45421-2024"
$ws.Range("B10").Value = "'45421"
$ws.Range("C10").Value = "'26865"
$ws.Range("E10").Value = $False
$ws.Range("F10").Value = "No screening recorded"
$ws.Range("G10").Value = $True
$ws.Range("H10").Value = $False

$ws.Range("A11").Value = "49143-2024"
$ws.Range("B11").Value = "'49143"
$ws.Range("C11").Value = "'38189"
$ws.Range("E11").Value = $False
$ws.Range("F11").Value = "No screening recorded"
$ws.Range("G11").Value = $True
$ws.Range("H11").Value = $True

$ws.Range("A12").Value = "56517-2024"
$ws.Range("B12").Value = "'56517"
$ws.Range("C12").Value = "'22777"
$ws.Range("E12").Value = $False
$ws.Range("F12").Value = "No screening recorded"
$ws.Range("G12").Value = $True
$ws.Range("H12").Value = $False

$ws.Range("A13").Value = "58898-2024"
$ws.Range("B13").Value = "'58898"
$ws.Range("C13").Value = "'93205"
$ws.Range("E13").Value = $False
$ws.Range("F13").Value = "No screening recorded"
$ws.Range("G13").Value = $False
$ws.Range("H13").Value = $False

$ws.Range("A14").Value = "64017-2024"
$ws.Range("B14").Value = "'64017"
$ws.Range("C14").Value = "'26431"
$ws.Range("E14").Value = $False
$ws.Range("F14").Value = "No screening recorded"
$ws.Range("G14").Value = $False
$ws.Range("H14").Value = $True

$ws.Range("A15").Value = "64607-2024"
$ws.Range("B15").Value = "'64607"
$ws.Range("C15").Value = "'11534"
$ws.Range("E15").Value = $False
$ws.Range("F15").Value = "No screening recorded"
$ws.Range("G15").Value = $True
$ws.Range("H15").Value = $False

$ws.Range("A16").Value = "76051-2024"
$ws.Range("B16").Value = "'76051"
$ws.Range("C16").Value = "'44751"
$ws.Range("E16").Value = $False
$ws.Range("F16").Value = "No screening recorded"
$ws.Range("G16").Value = $True
$ws.Range("H16").Value = $False

$ws.Range("A17").Value = "90185-2024"
$ws.Range("B17").Value = "'90185"
$ws.Range("C17").Value = "'37034"
$ws.Range("E17").Value = $False
$ws.Range("F17").Value = "No screening recorded"
$ws.Range("G17").Value = $True
$ws.Range("H17").Value = $True

$ws.Range("A18").Value = "92210-2024"
$ws.Range("B18").Value = "'92210"
$ws.Range("C18").Value = "'65158"
$ws.Range("E18").Value = $False
$ws.Range("F18").Value = "No screening recorded"
$ws.Range("G18").Value = $False
$ws.Range("H18").Value = $False

# --- D11 previously held a screening_encounter_id ("22777"); that value no longer applies ---
$ws.Range("D11").Value = ""
